# Generate Report for Handback
# This reflects a handback event: the localization status moves from
# "Ready for handoff" to "Handed back: in sync with en-US", the handback
# timestamps are refreshed, and the stale "handback not latest" error is
# cleared now that the content is in sync.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ----------------------------------------------------------------------
# Overview sheet - mirrors the per-language Status column
# ----------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E1").ColumnWidth = 29.14
$overview.Range("F1").ColumnWidth = 29.14

# ----------------------------------------------------------------------
# zh-cn sheet
# ----------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-10-20 00:25:39"
$zhcn.Range("P2").Value = ""
$zhcn.Range("C1").ColumnWidth = 29.14
$zhcn.Range("P1").ColumnWidth = 12.8

# ----------------------------------------------------------------------
# de-de sheet
# ----------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-10-20 00:25:57"
$dede.Range("P2").Value = ""
$dede.Range("C1").ColumnWidth = 29.14
$dede.Range("P1").ColumnWidth = 12.8
